$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "zcsm_1"
$ws.Range("A30").Value = "zcsm_3"
$ws.Range("A31").Value = "zcsm_6"
$ws.Range("A32").Value = "zcsm_12"
$ws.Range("A33").Value = "zcsm_24"
$ws.Range("A34").Value = "zcsm_36"

[void]$ws.Range("A35").Select()
